$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet1 (展览)
$ws1.Range("F4").Value = 4753
$ws1.Range("F5").Value = 204
$ws1.Range("F7").Value = 109
$ws1.Range("F10").Value = 0
$ws1.Range("F12").Value = 1120
$ws1.Range("F13").Value = 0
$ws1.Range("F14").Value = 259
$ws1.Range("F15").Value = 172
$ws1.Range("F17").Value = 139
$ws1.Range("F19").Value = 3840
$ws1.Range("F20").Value = 6172
$ws1.Range("F24").Value = 532
$ws1.Range("F25").Value = 0
$ws1.Range("F27").Value = 393
$ws1.Range("F28").Value = 0
$ws1.Range("F29").Value = 0
$ws1.Range("F30").Value = 567
$ws1.Range("F32").Value = 136
$ws1.Range("F33").Value = 0
$ws1.Range("F35").Value = 363
$ws1.Range("F36").Value = 159
$ws1.Range("F37").Value = 1553
$ws1.Range("F38").Value = 0
$ws1.Range("F41").Value = 55
$ws1.Range("F44").Value = 73
$ws1.Range("F45").Value = 0

# Sheet4 (全部类型)
$ws4.Range("F3").Value = 220
$ws4.Range("F5").Value = 204
$ws4.Range("F6").Value = 143
$ws4.Range("F7").Value = 109
$ws4.Range("F8").Value = 106
$ws4.Range("F13").Value = 0
$ws4.Range("F14").Value = 0
$ws4.Range("F16").Value = 0
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 139
$ws4.Range("F19").Value = 0
$ws4.Range("F21").Value = 0
$ws4.Range("F24").Value = 84
$ws4.Range("F25").Value = 532
$ws4.Range("F26").Value = 46
$ws4.Range("F27").Value = 3942
$ws4.Range("F28").Value = 0
$ws4.Range("F30").Value = 2538
$ws4.Range("F31").Value = 0
$ws4.Range("F33").Value = 0
$ws4.Range("F34").Value = 0
$ws4.Range("F35").Value = 292
$ws4.Range("F37").Value = 160
$ws4.Range("F38").Value = 1553
$ws4.Range("F39").Value = 0
$ws4.Range("F40").Value = 0
$ws4.Range("F41").Value = 51
$ws4.Range("F42").Value = 0
$ws4.Range("F43").Value = 0
$ws4.Range("F44").Value = 478
$ws4.Range("F45").Value = 0
$ws4.Range("F46").Value = 0

$wb.Save()
